# Updates cryptos list prices/volumes (commit: "Updated cryptos list ... with GitHub Actions")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "22.149.34"
$ws.Range("E2").Value = "  -1.47%  "

$ws.Range("D3").Value = "1.560.79"
$ws.Range("E3").Value = "  -1.01%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("E5").Value = "  -0.02%  "

$ws.Range("D6").Style = "Normal"
$ws.Range("D6").Value = "'289.94"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.43%  "

$ws.Range("D7").Style = "Normal"
$ws.Range("D7").Value = "'0.3806"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +3.34%  "

$ws.Range("D8").Style = "Normal"
$ws.Range("D8").Value = "'0.3285"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.48%  "

$ws.Range("D9").Style = "Normal"
$ws.Range("D9").Value = "'43.71"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -9.02%  "

$ws.Range("E10").Value = "  -0.84%  "

$ws.Range("E11").Value = "  -2.57%  "

$ws.Range("E12").Value = "  +0.00%  "

$ws.Range("D13").Style = "Normal"
$ws.Range("D13").Value = "'19.97"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.14%  "

$ws.Range("D14").Style = "Normal"
$ws.Range("D14").Value = "'5.829"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.26%  "

$ws.Range("D15").Style = "Normal"
$ws.Range("D15").Value = "'6.874"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.08%  "

$ws.Range("D16").Value = "1.563.56"
$ws.Range("E16").Value = "  -0.41%  "

$ws.Range("E17").Value = "  -2.57%  "

$ws.Range("D18").Style = "Normal"
$ws.Range("D18").Value = "'0.06638"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.39%  "

$ws.Range("D19").Style = "Normal"
$ws.Range("D19").Value = "'85.49"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.69%  "

$ws.Range("D20").Style = "Normal"
$ws.Range("D20").Value = "'6.460"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.89%  "

$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("E22").Value = "  -2.76%  "

$ws.Range("E23").Value = "  -2.35%  "

$ws.Range("D24").Value = "22.149.74"
$ws.Range("E24").Value = "  -1.39%  "

$ws.Range("D25").Style = "Normal"
$ws.Range("D25").Value = "'2.264"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -5.25%  "

$ws.Range("D26").Style = "Normal"
$ws.Range("D26").Value = "'2.539"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.99%  "

$ws.Range("D27").Style = "Normal"
$ws.Range("D27").Value = "'151.32"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.45%  "

$ws.Range("D28").Style = "Normal"
$ws.Range("D28").Value = "'19.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.05%  "

$ws.Range("D29").Style = "Normal"
$ws.Range("D29").Value = "'4.867"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.66%  "

$ws.Range("D30").Value = "1.739.21"
$ws.Range("E30").Value = "  -0.56%  "

$ws.Range("D31").Style = "Normal"
$ws.Range("D31").Value = "'121.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.00%  "

$ws.Range("D32").Style = "Normal"
$ws.Range("D32").Value = "'1.121"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.05%  "

$ws.Range("D33").Style = "Normal"
$ws.Range("D33").Value = "'6.029"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.99%  "

$ws.Range("D34").Style = "Normal"
$ws.Range("D34").Value = "'1.882"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -5.26%  "

$ws.Range("D35").Style = "Normal"
$ws.Range("D35").Value = "'9.348"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -5.39%  "

$ws.Range("D36").Style = "Normal"
$ws.Range("D36").Value = "'0.08214"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.91%  "

$ws.Range("D37").Style = "Normal"
$ws.Range("D37").Value = "'5.288"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.79%  "

$ws.Range("D38").Style = "Normal"
$ws.Range("D38").Value = "'0.02307"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.82%  "

$ws.Range("D39").Style = "Normal"
$ws.Range("D39").Value = "'0.06221"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.96%  "

$ws.Range("D40").Style = "Normal"
$ws.Range("D40").Value = "'0.2138"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.99%  "

$ws.Range("D41").Style = "Normal"
$ws.Range("D41").Value = "'1.230"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.91%  "

$ws.Range("D42").Style = "Normal"
$ws.Range("D42").Value = "'11.06"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.64%  "

$ws.Range("E43").Value = "  -0.01%  "

$ws.Range("D44").Style = "Normal"
$ws.Range("D44").Value = "'0.5980"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -5.19%  "

$ws.Range("D45").Style = "Normal"
$ws.Range("D45").Value = "'13.68"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.61%  "

$ws.Range("D46").Style = "Normal"
$ws.Range("D46").Value = "'3.761"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.77%  "

$ws.Range("E47").Value = "  -5.66%  "

$ws.Range("D48").Style = "Normal"
$ws.Range("D48").Value = "'1.989"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.87%  "

$ws.Range("D49").Style = "Normal"
$ws.Range("D49").Value = "'120.96"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.77%  "

$ws.Range("D50").Style = "Normal"
$ws.Range("D50").Value = "'1.171"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.57%  "

$ws.Range("D51").Style = "Normal"
$ws.Range("D51").Value = "'0.06991"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.32%  "
